$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value  = "成飞集成"
$ws.Range("C8").Value  = "领益智造"
$ws.Range("C9").Value  = "恒宝股份"
$ws.Range("C10").Value = "大位科技"
$ws.Range("C11").Value = "东华软件"
$ws.Range("C13").Value = "吉视传媒"
$ws.Range("C15").Value = "四川长虹"
$ws.Range("C17").Value = "奋达科技"
$ws.Range("C18").Value = "歌尔股份"
$ws.Range("C19").Value = "天融信"
$ws.Range("C20").Value = "大元泵业"
$ws.Range("C21").Value = "鸿博股份"
